$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new "department" column (C) for each course row.
$ws.Range("C2").Value = "Ageing Support"
$ws.Range("C3").Value = "Ageing Support"
$ws.Range("C4").Value = "Ageing Support"
$ws.Range("C5").Value = "Community Services"
$ws.Range("C6").Value = "Early Childhood"
$ws.Range("C7").Value = "Early Childhood"
$ws.Range("C8").Value = "Packages"
$ws.Range("C9").Value = "Packages"
$ws.Range("C10").Value = "Packages"
$ws.Range("C11").Value = "Packages"
$ws.Range("C12").Value = "Packages"

# Match the author's final selection/cursor position recorded in the sheet.
$ws.Range("C12").Select() | Out-Null
